# Add a new tracker row (#15 / row 16) documenting a new ADC-related device
# issue, matching the commit "Updated with ADCON2 = 15 ' Set Negative
# Reference Setting to ADNREF in ADCON1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newDescription = "New device with an ADC error.  Needs ADCON2 = 15 ' Set Negative Reference Setting to ADNREF in ADCON1 `nSee https://sourceforge.net/p/gcbasic/discussion/629990/thread/9b69d693/#e018"

# Index / Status / Description for the new row, mirroring the existing
# "OPEN" issue rows above it.
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "OPEN"
$ws.Cells.Item(16, 4).Value = $newDescription

# Match the wrapped, top/left aligned formatting used by the other
# Description cells (column D).
$ws.Cells.Item(16, 4).HorizontalAlignment = -4131
$ws.Cells.Item(16, 4).VerticalAlignment = -4160
$ws.Cells.Item(16, 4).WrapText = $true

# Same row height as the other multi-line entries of similar length.
$ws.Rows.Item(16).RowHeight = 45

# Reset the view back to the top of the sheet with D1 selected.
$ws.Range("D1").Select()
